# Weekly update: insert a new fruit/vegetable price record as row 47
# (Arándano (blue), "Primera" quality, fecha 2022-01-28, Provincia de
# Linares), pushing the existing rows 47-86 down to 48-87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 47 - shifts rows 47:86 down to 48:87
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the new price record
$ws.Cells.Item(47, 1).Value  = 11
$ws.Cells.Item(47, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(47, 3).Value  = "Bíobío"
$ws.Cells.Item(47, 4).Value  = 44589
$ws.Cells.Item(47, 5).Value  = 8
$ws.Cells.Item(47, 6).Value  = "Fruta"
$ws.Cells.Item(47, 7).Value  = 100101
$ws.Cells.Item(47, 8).Value  = "Berries"
$ws.Cells.Item(47, 9).Value  = 100101001
$ws.Cells.Item(47, 10).Value = "Arándano (blue)"
$ws.Cells.Item(47, 11).Value = "Sin especificar"
$ws.Cells.Item(47, 12).Value = "Primera"
$ws.Cells.Item(47, 13).Value = 220
$ws.Cells.Item(47, 14).Value = 3000
$ws.Cells.Item(47, 15).Value = 3500
$ws.Cells.Item(47, 16).Value = 3273
$ws.Cells.Item(47, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(47, 18).Value = "Provincia de Linares"
$ws.Cells.Item(47, 19).Value = 1636
$ws.Cells.Item(47, 20).Value = 2
